$d = $word.ActiveDocument

# Step 1: strip the trailing ". " (period + space) from the end of the
# paragraph, leaving "...supporting methods" (no trailing punctuation).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "as its supporting methods. ",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "as its supporting methods",
    2
)

# Step 2: append the new sentence as three additional runs, matching the
# author's edit (each insertion lands in its own run).
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertAfter(". At the time of writing, 5:00 PM 12/14/20")

$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertAfter(",")

$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertAfter(" Tom has not contributed to the project.")
